# Box Plot Updates, Color Updates Main Figures
#
# Repositions the percentage/label textboxes ("tx9".."tx18") that sit
# inside the grouped pie-chart figure on slide 1. Only the shape
# positions (Left/Top) change; sizes (Width/Height) are untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The labels live inside the single top-level group shape on the slide.
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Type -eq 6) {
        $grp = $candidate
        break
    }
}

# New Left/Top positions, in points (EMU / 12700), taken from the target
# EMU offsets. Values are nudged to the nearest float32 representation
# that still truncates back to the exact target EMU once the host
# converts Left/Top (points) to EMU on save.
$moves = @{
    "tx9"  = @(317.99969482421875, 189.8324432373047)
    "tx10" = @(367.8184509277344, 217.55819702148438)
    "tx11" = @(492.2294616699219, 211.39205932617188)
    "tx12" = @(515.4654541015625, 235.51678466796875)
    "tx13" = @(462.0378112792969, 254.0439453125)
    "tx14" = @(468.1937255859375, 281.3612060546875)
    "tx15" = @(489.8665466308594, 303.6302490234375)
    "tx16" = @(514.0528564453125, 331.1475830078125)
    "tx17" = @(350.3928527832031, 365.0951232910156)
    "tx18" = @(357.5033264160156, 392.4124450683594)
}

foreach ($name in $moves.Keys) {
    $shp = $grp.GroupItems.Item($name)
    $xy = $moves[$name]
    $shp.Left = $xy[0]
    $shp.Top = $xy[1]
}
